$p = $ppt.ActivePresentation

# Re-style the three tables (slides 14, 15, 16) from the plain custom
# "Table_0" style to the built-in table style, matching the
# Table Design ribbon action recorded in the source commit.
$newStyleId = "{3D389BA1-653B-4D7E-B891-3C729B42FE02}"

$tableSlideIndexes = @(14, 15, 16)
foreach ($idx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
